$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Days remaining" counters (one day has passed)
$ws.Range("B8").Value = 14
$ws.Range("B10").Value = 34

# The query refresh drops the old "TEST" / COLO-PREVENT trailing row and
# removes the stray cell style that had been applied to column A.
$ws.Range("A2:A11").Style = "Normal"
$ws.Range("C11").ClearContents()
$ws.Rows.Item(12).Delete()

# Keep the hidden ExternalData_1 defined name (and therefore the query
# table / autofilter range) in sync with the now-smaller data range.
$n = $wb.Names.Item("ExternalData_1")
$n.RefersTo = "=Sheet1!`$A`$1:`$C`$11"
